# Auto-generated edit script applying the diff to cryptos.xlsx
# Cell values are written with a leading apostrophe (forces the
# smart-entry parser to keep them as text instead of silently
# coercing numeric-looking strings into Number cells) and then
# ClearFormats() strips the resulting "quote prefix" style so the
# cell keeps the workbook-wide default style (no style churn).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'26.248.21"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.Value = "'  +1.78%  "
$c.ClearFormats()
$c = $ws.Range("D3")
$c.Value = "'1.607.25"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.Value = "'  +0.54%  "
$c.ClearFormats()
$c = $ws.Range("E4")
$c.Value = "'  -0.38%  "
$c.ClearFormats()
$c = $ws.Range("D5")
$c.Value = "'212.42"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.Value = "'  +1.76%  "
$c.ClearFormats()
$c = $ws.Range("E6")
$c.Value = "'  -0.44%  "
$c.ClearFormats()
$c = $ws.Range("D7")
$c.Value = "'0.483"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.Value = "'  +0.50%  "
$c.ClearFormats()
$c = $ws.Range("D8")
$c.Value = "'0.250"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.Value = "'  +1.31%  "
$c.ClearFormats()
$c = $ws.Range("E9")
$c.Value = "'  +1.39%  "
$c.ClearFormats()
$c = $ws.Range("D10")
$c.Value = "'18.24"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.Value = "'  +1.87%  "
$c.ClearFormats()
$c = $ws.Range("D11")
$c.Value = "'0.0799"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.Value = "'  +2.10%  "
$c.ClearFormats()
$c = $ws.Range("D12")
$c.Value = "'1.829.28"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.Value = "'  +0.37%  "
$c.ClearFormats()
$c = $ws.Range("D13")
$c.Value = "'1.605.82"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.Value = "'  +0.57%  "
$c.ClearFormats()
$c = $ws.Range("E14")
$c.Value = "'  -0.83%  "
$c.ClearFormats()
$c = $ws.Range("D15")
$c.Value = "'0.510"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.Value = "'  +0.19%  "
$c.ClearFormats()
$c = $ws.Range("D16")
$c.Value = "'26.205.56"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.Value = "'  +1.62%  "
$c.ClearFormats()
$c = $ws.Range("D17")
$c.Value = "'60.71"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.Value = "'  +0.42%  "
$c.ClearFormats()
$c = $ws.Range("D18")
$c.Value = "'0.0₃0729"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.Value = "'  +1.99%  "
$c.ClearFormats()
$c = $ws.Range("E19")
$c.Value = "'  -0.21%  "
$c.ClearFormats()
$c = $ws.Range("D20")
$c.Value = "'199.67"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.Value = "'  +5.50%  "
$c.ClearFormats()
$c = $ws.Range("D21")
$c.Value = "'4.25"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.Value = "'  +1.70%  "
$c.ClearFormats()
$c = $ws.Range("D22")
$c.Value = "'9.41"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.Value = "'  +0.82%  "
$c.ClearFormats()
$c = $ws.Range("D23")
$c.Value = "'6.01"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.Value = "'  +1.32%  "
$c.ClearFormats()
$c = $ws.Range("E24")
$c.Value = "'  +2.63%  "
$c.ClearFormats()
$c = $ws.Range("D25")
$c.Value = "'142.38"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.Value = "'  +1.13%  "
$c.ClearFormats()
$c = $ws.Range("E27")
$c.Value = "'  -0.35%  "
$c.ClearFormats()
$c = $ws.Range("D28")
$c.Value = "'15.19"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.Value = "'  +1.51%  "
$c.ClearFormats()
$c = $ws.Range("D29")
$c.Value = "'6.50"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.Value = "'  -0.25%  "
$c.ClearFormats()
$c = $ws.Range("E30")
$c.Value = "'  -0.82%  "
$c.ClearFormats()
$c = $ws.Range("E31")
$c.Value = "'  +0.85%  "
$c.ClearFormats()
$c = $ws.Range("E32")
$c.Value = "'  +2.20%  "
$c.ClearFormats()
$c = $ws.Range("E33")
$c.Value = "'  +0.41%  "
$c.ClearFormats()
$c = $ws.Range("E34")
$c.Value = "'  +2.03%  "
$c.ClearFormats()
$c = $ws.Range("D35")
$c.Value = "'2.35"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.Value = "'  -1.71%  "
$c.ClearFormats()
$c = $ws.Range("D36")
$c.Value = "'1.109.74"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.Value = "'  +1.37%  "
$c.ClearFormats()
$c = $ws.Range("E37")
$c.Value = "'  -0.61%  "
$c.ClearFormats()
$c = $ws.Range("B38")
$c.Value = "'VeChain"
$c.ClearFormats()
$c = $ws.Range("C38")
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.ClearFormats()
$c = $ws.Range("D38")
$c.Value = "'0.0152"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.Value = "'  +0.82%  "
$c.ClearFormats()
$c = $ws.Range("B39")
$c.Value = "'PaxDollar"
$c.ClearFormats()
$c = $ws.Range("C39")
$c.Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c.ClearFormats()
$c = $ws.Range("D39")
$c.Value = "'1.00"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.Value = "'  -0.02%  "
$c.ClearFormats()
$c = $ws.Range("B40")
$c.Value = "'ARBITRUM"
$c.ClearFormats()
$c = $ws.Range("C40")
$c.Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.Value = "'0.789"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.Value = "'  -0.17%  "
$c.ClearFormats()
$c = $ws.Range("B41")
$c.Value = "'ImmutableX"
$c.ClearFormats()
$c = $ws.Range("C41")
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.ClearFormats()
$c = $ws.Range("D41")
$c.Value = "'0.503"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.Value = "'  +1.18%  "
$c.ClearFormats()
$c = $ws.Range("D42")
$c.Value = "'0.782"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.Value = "'  +5.57%  "
$c.ClearFormats()
$c = $ws.Range("D43")
$c.Value = "'1.740.22"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.Value = "'  +0.29%  "
$c.ClearFormats()
$c = $ws.Range("E44")
$c.Value = "'  +1.17%  "
$c.ClearFormats()
$c = $ws.Range("D45")
$c.Value = "'92.90"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.Value = "'  -2.78%  "
$c.ClearFormats()
$c = $ws.Range("B46")
$c.Value = "'BabyDogeCoin"
$c.ClearFormats()
$c = $ws.Range("C46")
$c.Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c.ClearFormats()
$c = $ws.Range("D46")
$c.Value = "'0.0₆0108"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.Value = "'  -4.14%  "
$c.ClearFormats()
$c = $ws.Range("B47")
$c.Value = "'RenderToken"
$c.ClearFormats()
$c = $ws.Range("C47")
$c.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.Value = "'1.55"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.Value = "'  +9.64%  "
$c.ClearFormats()
$c = $ws.Range("D48")
$c.Value = "'53.67"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.Value = "'  +0.91%  "
$c.ClearFormats()
$c = $ws.Range("D49")
$c.Value = "'0.0510"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.Value = "'  -0.36%  "
$c.ClearFormats()
$c = $ws.Range("E50")
$c.Value = "'  -0.12%  "
$c.ClearFormats()
$c = $ws.Range("E51")
$c.Value = "'  -0.22%  "
$c.ClearFormats()
